$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the top of this price
# series. Insert a fresh row at 171 (pushing the existing rows 171-229
# down to 172-230) and fill it in with the new record.
$ws.Rows.Item(171).Insert()

$ws.Range("A171").Value = 4
$ws.Range("B171").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C171").Value = "Los Lagos"
$ws.Range("D171").Value = Get-Date -Year 2022 -Month 3 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Range("E171").Value = 10
$ws.Range("F171").Value = 100112003
$ws.Range("G171").Value = "Ajo"
$ws.Range("H171").Value = "Chino"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 80
$ws.Range("K171").Value = 21000
$ws.Range("L171").Value = 21000
$ws.Range("M171").Value = 21000
$ws.Range("N171").Value = "$/caja 10 kilos"
$ws.Range("O171").Value = "China"
$ws.Range("P171").Value = 2100
$ws.Range("Q171").Value = 10
$ws.Range("R171").Value = "Hortaliza"
